$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "CS604A[DC]  /  CS604B[SLa]"
$ws.Range("B4").Value = "CS601[SSK]  /  []"
$ws.Range("C4").Value = "CS605A[AH]  /  CS605B[SDe]"
$ws.Range("D4").Value = "Free Period!"
$ws.Range("E4").Value = "HU601[AnD]  /  []"
$ws.Range("F4").Value = "Free Period!"
$ws.Range("G4").Value = "CS603[BDu]  /  []"

$ws.Range("A6").Value = "CS604A[DC]  /  CS604B[SLa]"
$ws.Range("B6").Value = "CS601[SSK]  /  []"
$ws.Range("C6").Value = "Free Period!"
$ws.Range("D6").Value = "Free Period!"
$ws.Range("E6").Value = "CS602[DC]  /  []"
$ws.Range("F6").Value = "CS605A[AH]  /  CS605B[SDe]"
$ws.Range("G6").Value = "CS603[BDu]  /  []"

$ws.Range("A8").Value = "CS692[DC, GY]  /  CS691[SSK, SDe]"
$ws.Range("B8").Value = "CS692[DC, GY]  /  CS691[SSK, SDe]"
$ws.Range("C8").Value = "CS692[DC, GY]  /  CS691[SSK, SDe]"
$ws.Range("D8").Value = "Free Period!"
$ws.Range("E8").Value = "CS602[DC]  /  []"
$ws.Range("F8").Value = "Free Period!"
$ws.Range("G8").Value = "Free Period!"

$ws.Range("A10").Value = "CS605A[AH]  /  CS605B[SDe]"
$ws.Range("B10").Value = "CS601[SSK]  /  []"
$ws.Range("C10").Value = "Free Period!"
$ws.Range("D10").Value = "CS604A[DC]  /  CS604B[SLa]"
$ws.Range("E10").Value = "CS602[DC]  /  []"
$ws.Range("F10").Value = "Free Period!"
$ws.Range("G10").Value = "HU601[AnD]  /  []"

$ws.Range("A12").Value = "CS691[SSK, SDe]  /  CS693[BDu, AP]"
$ws.Range("B12").Value = "CS691[SSK, SDe]  /  CS693[BDu, AP]"
$ws.Range("C12").Value = "CS691[SSK, SDe]  /  CS693[BDu, AP]"
$ws.Range("D12").Value = "CS603[BDu]  /  []"
$ws.Range("E12").Value = "CS693[BDu, AP]  /  CS692[DC, SSK]"
$ws.Range("F12").Value = "CS693[BDu, AP]  /  CS692[DC, SSK]"
$ws.Range("G12").Value = "CS693[BDu, AP]  /  CS692[DC, SSK]"

$ws.Range("A16").Value = "Free Period!"
$ws.Range("B16").Value = "Free Period!"
$ws.Range("C16").Value = "Free Period!"
$ws.Range("D16").Value = "IT602[SU]  /  []"
$ws.Range("E16").Value = "IT603[ARC]  /  []"
$ws.Range("F16").Value = "IT604A[AB]  /  IT604B[RG]"
$ws.Range("G16").Value = "IT601A[AB]  /  IT601B[SL]"

$ws.Range("A18").Value = "IT693[SU, AD]  /  IT692[ARC, KDa]"
$ws.Range("B18").Value = "IT693[SU, AD]  /  IT692[ARC, KDa]"
$ws.Range("C18").Value = "IT693[SU, AD]  /  IT692[ARC, KDa]"
$ws.Range("D18").Value = "IT603[ARC]  /  []"
$ws.Range("E18").Value = "IT605[AGh]  /  []"
$ws.Range("F18").Value = "IT602[SU]  /  []"
$ws.Range("G18").Value = "HU601[SA]  /  []"

$ws.Range("A20").Value = "IT605[AGh]  /  []"
$ws.Range("B20").Value = "HU685[ACh, SSR]  /  []"
$ws.Range("C20").Value = "HU685[ACh, SSR]  /  []"
$ws.Range("D20").Value = "HU685[ACh, SSR]  /  []"
$ws.Range("E20").Value = "Free Period!"
$ws.Range("F20").Value = "Free Period!"
$ws.Range("G20").Value = "IT604A[AB]  /  IT604B[RG]"

$ws.Range("A22").Value = "IT692[ARC, KDa]  /  IT695[AGh, AB]"
$ws.Range("B22").Value = "IT692[ARC, KDa]  /  IT695[AGh, AB]"
$ws.Range("C22").Value = "IT692[ARC, KDa]  /  IT695[AGh, AB]"
$ws.Range("D22").Value = "IT604A[AB]  /  IT604B[RG]"
$ws.Range("E22").Value = "IT603[ARC]  /  []"
$ws.Range("F22").Value = "Free Period!"
$ws.Range("G22").Value = "IT601A[AB]  /  IT601B[SL]"

$ws.Range("A24").Value = "IT695[AGh, AB]  /  IT693[SU, AD]"
$ws.Range("B24").Value = "IT695[AGh, AB]  /  IT693[SU, AD]"
$ws.Range("C24").Value = "IT695[AGh, AB]  /  IT693[SU, AD]"
$ws.Range("D24").Value = "IT605[AGh]  /  []"
$ws.Range("E24").Value = "HU601[SA]  /  []"
$ws.Range("F24").Value = "IT602[SU]  /  []"
$ws.Range("G24").Value = "IT601A[AB]  /  IT601B[SL]"

$ws.Range("A28").Value = "ECE603A[PC]  /  ECE603B[JA]"
$ws.Range("B28").Value = "ECE604[AnC]  /  []"
$ws.Range("C28").Value = "HU601[TR]  /  []"
$ws.Range("D28").Value = "ECE601[SSG]  /  []"
$ws.Range("E28").Value = "ECE695[AD]  /  ECE692[SD]"
$ws.Range("F28").Value = "ECE695[AD]  /  ECE692[SD]"
$ws.Range("G28").Value = "ECE695[AD]  /  ECE692[SD]"

$ws.Range("A30").Value = "HU601[TR]  /  []"
$ws.Range("B30").Value = "ECE601[SSG]  /  []"
$ws.Range("C30").Value = "Free Period!"
$ws.Range("D30").Value = "Free Period!"
$ws.Range("E30").Value = "ECE603A[PC]  /  ECE603B[JA]"
$ws.Range("F30").Value = "ECE602[TD]  /  []"
$ws.Range("G30").Value = "ECE605A[AD]  /  []"

$ws.Range("A32").Value = "ECE604[AnC]  /  []"
$ws.Range("B32").Value = "Free Period!"
$ws.Range("C32").Value = "Free Period!"
$ws.Range("D32").Value = "ECE602[TD]  /  []"
$ws.Range("E32").Value = "ECE694[TD]  /  ECE695[AD]"
$ws.Range("F32").Value = "ECE694[TD]  /  ECE695[AD]"
$ws.Range("G32").Value = "ECE694[TD]  /  ECE695[AD]"

$ws.Range("A34").Value = "ECE603A[PC]  /  ECE603B[JA]"
$ws.Range("B34").Value = "ECE604[AnC]  /  []"
$ws.Range("C34").Value = "ECE605A[AD]  /  []"
$ws.Range("D34").Value = "ECE602[TD]  /  []"
$ws.Range("E34").Value = "ECE692[SD]  /  ECE694[TD]"
$ws.Range("F34").Value = "ECE692[SD]  /  ECE694[TD]"
$ws.Range("G34").Value = "ECE692[SD]  /  ECE694[TD]"

$ws.Range("A36").Value = "Free Period!"
$ws.Range("B36").Value = "HU685[ACh, AnC]  /  []"
$ws.Range("C36").Value = "HU685[ACh, AnC]  /  []"
$ws.Range("D36").Value = "HU685[ACh, AnC]  /  []"
$ws.Range("E36").Value = "Free Period!"
$ws.Range("F36").Value = "ECE601[SSG]  /  []"
$ws.Range("G36").Value = "ECE605A[AD]  /  []"

$ws.Range("A40").Value = "Free Period!"
$ws.Range("B40").Value = "EE602[ABo]  /  []"
$ws.Range("C40").Value = "EE604A[DC]  /  EE604B[GY]"
$ws.Range("D40").Value = "EE603[SKB]  /  []"
$ws.Range("E40").Value = "EE693[SDG, SMo]  /  EE694B[AP, DC]"
$ws.Range("F40").Value = "EE693[SDG, SMo]  /  EE694B[AP, DC]"
$ws.Range("G40").Value = "EE693[SDG, SMo]  /  EE694B[AP, DC]"

$ws.Range("A42").Value = "EE605A[IB]  /  EE605B[JA]"
$ws.Range("B42").Value = "EE694B[GY, SDe]  /  EE693[SDG, SMo]"
$ws.Range("C42").Value = "EE694B[GY, SDe]  /  EE693[SDG, SMo]"
$ws.Range("D42").Value = "EE694B[GY, SDe]  /  EE693[SDG, SMo]"
$ws.Range("E42").Value = "Free Period!"
$ws.Range("F42").Value = "EE603(T)[SKB, SDG]  /  []"
$ws.Range("G42").Value = "EE603[SKB]  /  []"

$ws.Range("A44").Value = "EE691[PG, ASG]  /  EE692[ABo, SDC]"
$ws.Range("B44").Value = "EE691[PG, ASG]  /  EE692[ABo, SDC]"
$ws.Range("C44").Value = "EE691[PG, ASG]  /  EE692[ABo, SDC]"
$ws.Range("D44").Value = "EE601[PG]  /  []"
$ws.Range("E44").Value = "EE605A[IB]  /  EE605B[SD]"
$ws.Range("F44").Value = "EE604A[DC]  /  EE604B[GY]"
$ws.Range("G44").Value = "EE602[ABo]  /  []"

$ws.Range("A46").Value = "EE601[PG]  /  []"
$ws.Range("B46").Value = "EE602(T)[ABo, SDC]  /  []"
$ws.Range("C46").Value = "Free Period!"
$ws.Range("D46").Value = "EE601(T)[PG, ASG]  /  []"
$ws.Range("E46").Value = "EE692[ABo, SDC]  /  EE691[PG, ASG]"
$ws.Range("F46").Value = "EE692[ABo, SDC]  /  EE691[PG, ASG]"
$ws.Range("G46").Value = "EE692[ABo, SDC]  /  EE691[PG, ASG]"

$ws.Range("A48").Value = "EE605A[IB]  /  EE605B[JA]"
$ws.Range("B48").Value = "EE602[ABo]  /  []"
$ws.Range("C48").Value = "EE604A[DC]  /  EE604B[GY]"
$ws.Range("D48").Value = "EE601[PG]  /  []"
$ws.Range("E48").Value = "Free Period!"
$ws.Range("F48").Value = "Free Period!"
$ws.Range("G48").Value = "EE603[SKB]  /  []"
